# Add a new doctor schedule row for "Dr. rubin" on Friday, 8:00 AM - 7:00 PM
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("sheet2")

$ws.Range("A11:B11").Copy()
$ws.Range("A12:B12").PasteSpecial(-4122)

$ws.Cells.Item(12, 1).Value = "Dr. rubin"
$ws.Cells.Item(12, 2).Value = "Friday"
$ws.Cells.Item(12, 3).Value = 0.33333333333333331
$ws.Cells.Item(12, 4).Value = 0.79166666666666663

$ws.Range("C12:D12").NumberFormat = "h:mm"
